# "31 Jan Presenti Sheet" — add a new attendance column (AH) for 31-Jan-2024
# on the Jan-2024 sheet, mirroring the formatting of the previous day's
# column (AG), and extend the data-validation list range to cover it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the whole AG1:AG4 column (values + formatting) into AH1:AH4 so the
# new column inherits the exact same styles (date header style, border,
# shared-string "Present"/"Absent" entries, etc.) as the rest of the table.
$ws.Range("AG1:AG4").Copy($ws.Range("AH1:AH4"))

# AH1 is the new date header — 31-Jan-2024 (Excel serial 45322), one day
# after AG1 (45321).
$ws.Range("AH1").Value = 45322

# Extend the "Present, Absent,Reason" list validation so it also governs
# the newly added AH2:AH4 cells (was C2:AG4, now C2:AH4).
$ws.Range("C2:AG4").Validation.Delete()
$ws.Range("C2:AH4").Validation.Add(3, 1, 1, '"Present, Absent,Reason"')

# Matches the author's last selection in the saved file.
$ws.Range("AE18").Select()
